$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Test expr" sample-call table in columns L:M (rows 13-17) ---
$ws.Range("L13").Value = "Test expr"

# --- Update the existing "tte" sample-call table header (H14:I14 merged) ---
# The method under test changed from tte(boolean test) to expr(boolean in).
$ws.Range("H14").Value = "SimpleRules Foo[] expr(boolean in)"

$ws.Range("L14").Value = "in"
$ws.Range("M14").Value = "_res_"

$ws.Range("L15").Value = "Input"
$ws.Range("M15").Value = "Result"

$ws.Range("L16").Value = $true
$ws.Range("M16").Value = "1,3,5"

$ws.Range("L17").Value = $false
$ws.Range("M17").Value = "2,4,6"

# Move the active selection to reflect where editing ended up.
$ws.Range("M18").Select() | Out-Null
